# Update spritesheet reference table
# - Player "Fire" row becomes "Attack" (now repeats) and prefix changes
# - Enemy 01 "Fire" row becomes "Move" (duplicate Move entry, prefix changes)
# - Enemy 02 "Fire" row becomes "Attack" (now repeats) and prefix changes
# - Enemy 03 rows shrink from 48x48 to 36x36, and gain Attack / Aim rows
#   in place of the old Fire / (duplicate) Fire rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sprites")

# Row 6: Player / Fire -> Attack
$ws.Cells.Item(6, 2).Value = "Attack"
$ws.Cells.Item(6, 5).Value = "yes"
$ws.Cells.Item(6, 9).Value = "player/attack-##"

# Row 11: Enemy 01 / Fire -> Move
$ws.Cells.Item(11, 2).Value = "Move"
$ws.Cells.Item(11, 5).Value = "yes"
$ws.Cells.Item(11, 9).Value = "enemy01/move-##"

# Row 14: Enemy 02 / Fire -> Attack
$ws.Cells.Item(14, 2).Value = "Attack"
$ws.Cells.Item(14, 5).Value = "yes"
$ws.Cells.Item(14, 9).Value = "enemy02/attack-##"

# Row 16: Enemy 03 / Move - resize 48 -> 36
$ws.Cells.Item(16, 6).Value = 36
$ws.Cells.Item(16, 7).Value = 36

# Row 17: Enemy 03 / Fire -> Attack, resize 48 -> 36
$ws.Cells.Item(17, 2).Value = "Attack"
$ws.Cells.Item(17, 6).Value = 36
$ws.Cells.Item(17, 7).Value = 36
$ws.Cells.Item(17, 9).Value = "enemy03/attack-##"

# Row 18: Enemy 03 / Hit -> Aim, resize 48 -> 36
$ws.Cells.Item(18, 2).Value = "Aim"
$ws.Cells.Item(18, 6).Value = 36
$ws.Cells.Item(18, 7).Value = 36
$ws.Cells.Item(18, 9).Value = "enemy03/aim-##"

# Row 19: Enemy 03 / Die - resize 48 -> 36
$ws.Cells.Item(19, 6).Value = 36
$ws.Cells.Item(19, 7).Value = 36

$ws.Range("E18").Select()
